# Reorders the ATD class rows in the rats subgroups table (Plan1 sheet).
# The row that held the "MAOI / selegiline" data (row 93) is moved up to
# row 90, and the three rows that used to be 90-92 ("All TeCA" entries)
# each shift down by one row (90->91, 91->92, 92->93).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C","D","E","F","G","H","I","J","K")

# Snapshot the current values for columns C:K across rows 90-93 before
# making any changes, so that later writes do not clobber data we still
# need to read.
$snapshot = @{}
foreach ($r in 90..93) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Destination row -> source row (old row 93 becomes the new row 90, and
# the old 90/91/92 rows shift down by one).
$mapping = @{ 90 = 93; 91 = 90; 92 = 91; 93 = 92 }

foreach ($destRow in 90..93) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $newVal = $srcVals[$c]
        $oldVal = $snapshot[$destRow][$c]
        # Only write the cell if the value actually changes, to avoid
        # needlessly re-serializing values that are already correct.
        if ($newVal -ne $oldVal) {
            $ws.Range("$c$destRow").Value = $newVal
        }
    }
}

# Update the sheet view to match the new scroll position / selection.
$ws.Range("C87").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 79
$win.ScrollColumn = 1
